$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 111, shifting existing rows 111-160 down to 112-161
$ws.Rows.Item(111).Insert()

# Populate the new row 111 with values (copy pattern from old row, with updated D/J/K/L/M/P)
$ws.Cells.Item(111, 1).Value = 8
$ws.Cells.Item(111, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(111, 3).Value = "Coquimbo"
$ws.Cells.Item(111, 4).Value = 44784
$ws.Cells.Item(111, 5).Value = 4
$ws.Cells.Item(111, 6).Value = 100112044
$ws.Cells.Item(111, 7).Value = "Perejil"
$ws.Cells.Item(111, 8).Value = "Sin especificar"
$ws.Cells.Item(111, 9).Value = "Primera"
$ws.Cells.Item(111, 10).Value = 2600
$ws.Cells.Item(111, 11).Value = 2000
$ws.Cells.Item(111, 12).Value = 2500
$ws.Cells.Item(111, 13).Value = 2250
$ws.Cells.Item(111, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(111, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(111, 16).Value = 1500
$ws.Cells.Item(111, 17).Value = 1.5
$ws.Cells.Item(111, 18).Value = "Hortaliza"
